# LnCO as on Dec 22
# Insert two new "weightage" key/value rows into the "Main" config sheet:
#   - new row 11: Output_Comparatives_Weightage_sheetname / Comparatives_top_weightage
#   - new row 18 (after the Concentration block shifts down): Output_Concentration_Weightage_sheetname / Concentration_top_weightage
# Everything below each insertion point shifts down accordingly. Fix up the
# hyperlinks (they don't follow the shifted cells automatically) and the
# data-validation formula (its range shifts, but the Formula1 text does not).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Insert new row 11: Comparatives "weightage" entry --------------------
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Output_Comparatives_Weightage_sheetname"
$ws.Range("B11").Value = "Comparatives_top_weightage"
$ws.Range("A11:B11").Style = "Normal"

# --- Insert new row 18: Concentration "weightage" entry -------------------
# (Before this insert, the old row 17 "Output_Concentration_Vendor_sheetname"
#  has already become row 17 due to the previous insert; the new row is
#  added right after it.)
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "Output_Concentration_Weightage_sheetname"
$ws.Range("B18").Value = "Concentration_top_weightage"
$ws.Range("A18:B18").Style = "Normal"

# --- Fix the data validation list on the env-file cells --------------------
# The validated range itself (B35:B36 -> B37:B38) shifts automatically with
# the row inserts, but the Formula1 text referencing $A$33:$A$34 needs to be
# repointed to the new LOCAL/QUALITY_DOTENV_FILE key rows ($A$35:$A$36).
$dv = $ws.Range("B37:B38").Validation
$dv.Formula1 = '$A$35:$A$36'

# --- Fix the mail-address hyperlinks ---------------------------------------
# Hyperlink anchors stay pinned to their original row/col through row
# inserts instead of following the moved cell content, so re-create them at
# the new locations (To_Mail_Address / CC_Mail_Address moved from B43/B44 to
# B45/B46).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B45"), "mailto:kalyan.gundu@bradsol.com")
$ws.Hyperlinks.Add($ws.Range("B46"), "mailto:kalyan.gundu@bradsol.com")
# Adding a hyperlink stamps a fresh (duplicate) "Hyperlink" style record on the
# cell instead of reusing the existing one already used elsewhere in the
# sheet; re-applying the named style snaps it back onto the pre-existing xf.
$ws.Range("B45").Style = "Hyperlink"
$ws.Range("B46").Style = "Hyperlink"

# --- Update the saved selection/view ---------------------------------------
$ws.Range("B8").Select()
